$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.172505855560303
$ws.Range("B1").Value = 5.942238807678223
$ws.Range("C1").Value = 3.991207838058472
$ws.Range("D1").Value = 1.761365413665771
$ws.Range("E1").Value = 1.212256073951721
